$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 2443.5908
$ws.Range("I53").Value = 623.875
$ws.Range("J53").Value = 3483.4285
$ws.Range("K53").Value = 623.875
$ws.Range("L53").Value = 3483.4285
$ws.Range("M53").Value = 13.125
$ws.Range("N53").Value = -4757.4285
$ws.Range("H64").Value = 12501.5
$ws.Range("H67").Value = 12501.5
$ws.Range("H76").Value = 6909.56
$ws.Range("J76").Value = 8813
$ws.Range("L76").Value = 8813
$ws.Range("N76").Value = -9443
$ws.Range("H79").Value = 6909.56
$ws.Range("J79").Value = 8813
$ws.Range("L79").Value = 8813
$ws.Range("N79").Value = -10997
$ws.Range("H98").Value = 464546.06
$ws.Range("I98").Value = 2544.6667
$ws.Range("J98").Value = 1504049.2
$ws.Range("K98").Value = 2544.6667
$ws.Range("L98").Value = 1504049.2
$ws.Range("M98").Value = -1046.6667
$ws.Range("N98").Value = -1507045.2
$ws.Range("H101").Value = 1440
$ws.Range("I101").Value = 1260.3077
$ws.Range("J101").Value = 2218.6667
$ws.Range("K101").Value = 3780.9231
$ws.Range("L101").Value = 6656.000100000001
$ws.Range("M101").Value = -2158.9231
$ws.Range("N101").Value = -9900.000100000001
$ws.Range("H103").Value = 1158.2307
$ws.Range("I103").Value = 712.6667
$ws.Range("J103").Value = 1291.9
$ws.Range("K103").Value = 2138.0001
$ws.Range("L103").Value = 3875.7
$ws.Range("M103").Value = -1552.0001
$ws.Range("N103").Value = -5047.700000000001
$ws.Range("H112").Value = 1944.5428
$ws.Range("J112").Value = 1944.5428
$ws.Range("L112").Value = 5833.6284
$ws.Range("N112").Value = -8049.6284
$ws.Range("H122").Value = 464546.06
$ws.Range("I122").Value = 2544.6667
$ws.Range("J122").Value = 1504049.2
$ws.Range("K122").Value = 7634.000100000001
$ws.Range("L122").Value = 4512147.6
$ws.Range("M122").Value = -5184.000100000001
$ws.Range("N122").Value = -4517047.6
$ws.Range("H135").Value = 6221.4614
$ws.Range("I135").Value = 5538.9
$ws.Range("K135").Value = 49850.1
$ws.Range("M135").Value = -47315.1
$ws.Range("H140").Value = 61899
$ws.Range("J140").Value = 61899
$ws.Range("L140").Value = 61899
$ws.Range("N140").Value = -72259

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15956.154
$ws.Range("I32").Value = 10211.5
$ws.Range("K32").Value = 10211.5
$ws.Range("M32").Value = -9924.5
$ws.Range("H56").Value = 50000
$ws.Range("J56").Value = 50000
$ws.Range("L56").Value = 50000
$ws.Range("N56").Value = -51484
$ws.Range("H61").Value = 6437.857
$ws.Range("I61").Value = 5454.6113
$ws.Range("K61").Value = 5454.6113
$ws.Range("M61").Value = -5242.6113
$ws.Range("H63").Value = 7937.5
$ws.Range("I63").Value = 2750
$ws.Range("K63").Value = 2750
$ws.Range("M63").Value = -2064
$ws.Range("H66").Value = 7937.5
$ws.Range("I66").Value = 2750
$ws.Range("K66").Value = 13750
$ws.Range("M66").Value = -10318
$ws.Range("H132").Value = 5311.375
$ws.Range("I132").Value = 4711
$ws.Range("K132").Value = 14133
$ws.Range("M132").Value = -11603
$ws.Range("H136").Value = 6437.857
$ws.Range("I136").Value = 5454.6113
$ws.Range("K136").Value = 16363.8339
$ws.Range("M136").Value = -13813.8339

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2575.0435
$ws.Range("I99").Value = 1682.4736
$ws.Range("K99").Value = 1682.4736
$ws.Range("M99").Value = -184.4736
$ws.Range("H105").Value = 22998.4
$ws.Range("I105").Value = 34283.168
$ws.Range("K105").Value = 34283.168
$ws.Range("M105").Value = -32536.168
$ws.Range("H134").Value = 4858.3
$ws.Range("I134").Value = 4456
$ws.Range("J134").Value = 7138
$ws.Range("K134").Value = 13368
$ws.Range("L134").Value = 21414
$ws.Range("M134").Value = -10833
$ws.Range("N134").Value = -26484

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37096.312
$ws.Range("I31").Value = 3437.7778
$ws.Range("K31").Value = 3437.7778
$ws.Range("M31").Value = -3142.7778
$ws.Range("H34").Value = 37096.312
$ws.Range("I34").Value = 3437.7778
$ws.Range("K34").Value = 3437.7778
$ws.Range("M34").Value = -3235.7778
$ws.Range("H62").Value = 9674.6
$ws.Range("J62").Value = 29006
$ws.Range("L62").Value = 29006
$ws.Range("N62").Value = -30254
$ws.Range("H65").Value = 9674.6
$ws.Range("J65").Value = 29006
$ws.Range("L65").Value = 145030
$ws.Range("N65").Value = -151270
$ws.Range("H105").Value = 2703.5
$ws.Range("J105").Value = 3365.4443
$ws.Range("L105").Value = 3365.4443
$ws.Range("N105").Value = -6859.4443

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 54628.727
$ws.Range("J2").Value = 85805.28999999999
$ws.Range("L2").Value = 514831.74
$ws.Range("N2").Value = -515057.74
$ws.Range("H122").Value = 4359
$ws.Range("J122").Value = 4516.2666
$ws.Range("L122").Value = 40646.3994
$ws.Range("N122").Value = -45546.3994
$ws.Range("H132").Value = 5444.5386
$ws.Range("I132").Value = 4080
$ws.Range("K132").Value = 36720
$ws.Range("M132").Value = -34190

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1992.3125
$ws.Range("I102").Value = 642.4545000000001
$ws.Range("J102").Value = 4962
$ws.Range("K102").Value = 642.4545000000001
$ws.Range("L102").Value = 4962
$ws.Range("M102").Value = 979.5454999999999
$ws.Range("N102").Value = -8206
$ws.Range("H113").Value = 4195.923
$ws.Range("J113").Value = 4283.3335
$ws.Range("L113").Value = 4283.3335
$ws.Range("N113").Value = -8623.333500000001
$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 15000
$ws.Range("M132").Value = -12470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 9038.538
$ws.Range("I68").Value = 5499.3335
$ws.Range("K68").Value = 5499.3335
$ws.Range("M68").Value = -4750.3335
$ws.Range("H71").Value = 9038.538
$ws.Range("I71").Value = 5499.3335
$ws.Range("K71").Value = 27496.6675
$ws.Range("M71").Value = -23752.6675
$ws.Range("H122").Value = 197939.95
$ws.Range("I122").Value = 273615
$ws.Range("J122").Value = 8752.333000000001
$ws.Range("K122").Value = 820845
$ws.Range("L122").Value = 26256.999
$ws.Range("M122").Value = -818395
$ws.Range("N122").Value = -31156.999
$ws.Range("H132").Value = 8661.923000000001
$ws.Range("I132").Value = 6525
$ws.Range("J132").Value = 12081
$ws.Range("K132").Value = 19575
$ws.Range("L132").Value = 36243
$ws.Range("M132").Value = -17045
$ws.Range("N132").Value = -41303

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6648.0835
$ws.Range("I126").Value = 2955.6
$ws.Range("K126").Value = 8866.799999999999
$ws.Range("M126").Value = -6396.799999999999
